# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# These two sheets mirror the same underlying data, so the same logical
# rows (by event) need to be bumped, though the row numbers differ
# slightly between the two sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    4  = 10351
    5  = 768
    8  = 490
    10 = 483
    11 = 281
    13 = 12942
    14 = 12942
    29 = 2181
    30 = 1159
    31 = 4350
    33 = 3936
    34 = 967
    35 = 2687
    36 = 3113
    38 = 1412
    40 = 800
    41 = 61
    42 = 170
    43 = 622
    44 = 882
    46 = 177
    47 = 368
    48 = 130
    49 = 203
    50 = 238
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    4  = 10351
    5  = 768
    7  = 490
    9  = 483
    10 = 281
    11 = 12942
    12 = 12942
    26 = 2181
    27 = 1159
    31 = 4350
    32 = 3936
    33 = 967
    34 = 2687
    35 = 3113
    40 = 800
    41 = 61
    42 = 622
    44 = 882
    46 = 177
    47 = 368
    48 = 130
    49 = 203
    50 = 238
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
